# Update "想去人数" (interested-count) values in column F across the four
# worksheets, reflecting the newly generated site stats (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1507
$ws1.Range("F5").Value  = 7656
$ws1.Range("F6").Value  = 4838
$ws1.Range("F22").Value = 1197
$ws1.Range("F35").Value = 118
$ws1.Range("F41").Value = 95

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F27").Value = 642
$ws2.Range("F28").Value = 49

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 70
$ws3.Range("F9").Value  = 75
$ws3.Range("F10").Value = 1672
$ws3.Range("F11").Value = 2574

# Sheet "全部类型" (All Types) -- aggregated view of the above sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1507
$ws4.Range("F9").Value  = 7656
$ws4.Range("F10").Value = 4838
$ws4.Range("F18").Value = 1672
$ws4.Range("F19").Value = 2574
$ws4.Range("F26").Value = 1197
$ws4.Range("F27").Value = 642
$ws4.Range("F35").Value = 118
$ws4.Range("F41").Value = 95
